$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at positions 9 and 10, pushing the existing
# rows 9-22 down to 11-24 (and expanding the sheet dimension to A1:T24).
$ws.Rows("9:10").Insert()

# --- New row 9: Damasco, Castle Brite, Primera, Región de O'Higgins ---
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "Macroferia Regional de Talca"
$ws.Range("C9").Value = "Maule"
$ws.Range("D9").Value = 44540
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100103
$ws.Range("H9").Value = "Frutos de hueso (carozo)"
$ws.Range("I9").Value = 100103003
$ws.Range("J9").Value = "Damasco"
$ws.Range("K9").Value = "Castle Brite"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("Q9").Value = "`$/caja 16 kilos"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1250
$ws.Range("T9").Value = 16

# --- New row 10: Damasco, Castle Brite, Segunda, Región de O'Higgins ---
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Macroferia Regional de Talca"
$ws.Range("C10").Value = "Maule"
$ws.Range("D10").Value = 44540
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103003
$ws.Range("J10").Value = "Damasco"
$ws.Range("K10").Value = "Castle Brite"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("Q10").Value = "`$/caja 16 kilos"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 938
$ws.Range("T10").Value = 16
